$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.281384468078613
$ws.Range("B1").Value = 1.650769591331482
$ws.Range("C1").Value = 3.030357599258423
$ws.Range("D1").Value = 1.505936741828918
$ws.Range("E1").Value = 0.8240828514099121
